$wb = $excel.ActiveWorkbook

# --- "Iteration #2" sheet: update hours logged and self-evaluation ---
$iter2 = $wb.Worksheets.Item("Iteration #2")
$iter2.Range("C14").Value = 4
$iter2.Range("C18").Value = 4
$iter2.Range("B40").Value = 7
$iter2.Range("B42").Value = "Aucun commentaire, tout c'est bien déroulé."

# --- "Iteration #3" sheet: log the first SQLite-implementation entry ---
$iter3 = $wb.Worksheets.Item("Iteration #3")
$iter3.Range("A14").Value = "4/23/2018"
$iter3.Range("B14").Value = "Debut de l'implémentation de SQLite dans l'app"
$iter3.Range("C14").Value = 3

# --- Make "Iteration #3" the active / selected sheet & cell, matching the
#     workbook's new activeTab and the sheet's new selection. ---
$iter3.Activate() | Out-Null
$iter3.Range("C14").Select() | Out-Null
